# Apply updated Seasonality Index / MyForecast values on "Forecast Comparison"
# and the refreshed 16-week forecast total on "Summary".

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet: column L = Seasonality Index, column D = MyForecast ---

# Column L (Seasonality Index) updates, rows 2-17
$wsForecast.Range("L2").Value  = 0.97
$wsForecast.Range("L3").Value  = 0.91
$wsForecast.Range("L4").Value  = 1.02
$wsForecast.Range("L5").Value  = 1.12
$wsForecast.Range("L6").Value  = 0.86
$wsForecast.Range("L7").Value  = 1
$wsForecast.Range("L8").Value  = 1.15
$wsForecast.Range("L9").Value  = 0.9399999999999999
$wsForecast.Range("L10").Value = 0.95
$wsForecast.Range("L11").Value = 0.96
$wsForecast.Range("L12").Value = 1.13
$wsForecast.Range("L13").Value = 0.83
$wsForecast.Range("L14").Value = 0.99
$wsForecast.Range("L15").Value = 0.89
$wsForecast.Range("L16").Value = 1.03
$wsForecast.Range("L17").Value = 1.11

# Column D (MyForecast) updates
$wsForecast.Range("D11").Value = 25
$wsForecast.Range("D12").Value = 25
$wsForecast.Range("D13").Value = 26
$wsForecast.Range("D15").Value = 25

# --- "Summary" sheet: Total Forecast (16 Weeks) ---
# Force text storage (matches existing inline-string "number-as-text" cells in this column)
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "414"
